$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data updates scraped for cryptos list (GitHub Actions refresh)

# Row 2
$ws.Range("D2").Value = "63.988.06"
$ws.Range("E2").Value = "  +1.41%  "

# Row 3
$ws.Range("D3").Value = "3.316.00"
$ws.Range("E3").Value = "  +6.03%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.96"
$ws.Range("E5").Value = "  +1.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.51"
$ws.Range("E6").Value = "  +5.16%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").Value = "3.316.40"
$ws.Range("E8").Value = "  +6.26%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  +1.47%  "

# Row 10
$ws.Range("E10").Value = "  +3.26%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.56"
$ws.Range("E11").Value = "  +6.20%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.474"
$ws.Range("E12").Value = "  +4.18%  "

# Row 13
$ws.Range("E13").Value = "  +1.71%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.98"
$ws.Range("E14").Value = "  +2.45%  "

# Row 15
$ws.Range("D15").Value = "3.868.41"
$ws.Range("E15").Value = "  +6.21%  "

# Row 16
$ws.Range("E16").Value = "  +0.02%  "

# Row 17
$ws.Range("D17").Value = "3.320.46"
$ws.Range("E17").Value = "  +6.18%  "

# Row 18
$ws.Range("D18").Value = "64.075.57"
$ws.Range("E18").Value = "  +1.53%  "

# Row 19
$ws.Range("E19").Value = "  +3.67%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "482.73"
$ws.Range("E20").Value = "  +2.26%  "

# Row 21
$ws.Range("E21").Value = "  +1.35%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.739"
$ws.Range("E22").Value = "  +6.11%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.01"
$ws.Range("E23").Value = "  +4.21%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.70"
$ws.Range("E24").Value = "  +5.79%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.02"
$ws.Range("E25").Value = "  -1.80%  "

# Row 26
$ws.Range("E26").Value = "  +0.22%  "

# Row 27
$ws.Range("E27").Value = "  +2.40%  "

# Row 28
$ws.Range("B28").Value = "NEARProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.30"
$ws.Range("E28").Value = "  +2.52%  "

# Row 29
$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.05%  "

# Row 30
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.26"
$ws.Range("E30").Value = "  +4.04%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "29.69"
$ws.Range("E31").Value = "  +11.16%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.17"
$ws.Range("E32").Value = "  +5.88%  "

# Row 33
$ws.Range("E33").Value = "  -1.97%  "

# Row 34
$ws.Range("E34").Value = "  +2.35%  "

# Row 35
$ws.Range("E35").Value = "  +2.60%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.01"
$ws.Range("E36").Value = "  +3.59%  "

# Row 37
$ws.Range("D37").Value = "0.0₃0761"
$ws.Range("E37").Value = "  +7.21%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "53.34"
$ws.Range("E38").Value = "  +2.52%  "

# Row 39
$ws.Range("E39").Value = "  +4.55%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "436.02"
$ws.Range("E40").Value = "  +2.79%  "

# Row 41
$ws.Range("D41").Value = "3.069.52"
$ws.Range("E41").Value = "  +5.98%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.79"
$ws.Range("E42").Value = "  +3.57%  "

# Row 43
$ws.Range("E43").Value = "  +2.86%  "

# Row 44
$ws.Range("E44").Value = "  -0.57%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.268"
$ws.Range("E45").Value = "  +2.61%  "

# Row 46
$ws.Range("E46").Value = "  +4.70%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.62"
$ws.Range("E47").Value = "  +4.07%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "36.04"
$ws.Range("E48").Value = "  +14.52%  "

# Row 50
$ws.Range("E50").Value = "  +2.96%  "

# Row 51
$ws.Range("E51").Value = "  +1.88%  "
